# geometryVariables.xlsx -- "minor modifications" commit
#
# Changes applied on the "Data" sheet:
#   1. A new geometry-variable row is inserted before the old row 27
#      ("x_ac"), which pushes that row and everything below it down by
#      one (old rows 27-37 become 28-38). The new row holds:
#         B: x_le_w   C: 19 (highlighted)   D: m
#         E: Distance of leading edge of wing from nose
#   2. C3  (x_frontSpar) value 15 -> 20.2, and gets the same yellow
#      "filled in" highlight used by the other populated value cells.
#   3. C4  (x_rearSpar)  value 18 -> 24.3, same highlight added.
#   4. C9  (wingBoxLoc)  value 18 -> 24.3 (no style change).
#   5. The fuelTankLen row (old C33, now C34 after the insert) had a
#      stray "no-fill" style applied -- it is reset to the normal
#      yellow highlight used by the rest of the sheet.
#   6. The active selection is moved to C9 (matching the saved view).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- 1. Insert the new "x_le_w" row above the old row 27 ---------------
[void]$ws.Rows.Item(27).Insert()

$ws.Range("E27").Value = "Distance of leading edge of wing from nose"
$ws.Range("B27").Value = "x_le_w"
$ws.Range("C27").Value = 19
$ws.Range("C27").Interior.Color = 65535
$ws.Range("D27").Value = "m"

# --- 2-4. Updated geometry values ---------------------------------------
$ws.Range("C3").Value = 20.2
$ws.Range("C3").Interior.Color = 65535

$ws.Range("C4").Value = 24.3
$ws.Range("C4").Interior.Color = 65535

$ws.Range("C9").Value = 24.3

# --- 5. Fix fuelTankLen's highlight (shifted down to row 34) -----------
$ws.Range("C34").Interior.Color = 65535

# --- 6. Restore the saved selection -------------------------------------
[void]$ws.Range("C9").Select()

Write-Output "Applied geometryVariables edits"
